# Update the "想去人数" (number of people interested) figures (column F)
# for several events on the "展览", "演出" and "全部类型" sheets, reflecting
# the latest scrape of the source data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet (sheetId 1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 131    # 苏州·无限次元夜场
$wsExhibit.Range("F8").Value = 141    # 苏州·漫语堂动漫嘉年华
$wsExhibit.Range("F9").Value = 330    # 苏州·第三届华盟国漫次元嘉年华
$wsExhibit.Range("F13").Value = 11578 # 苏州·COME IN JOY 动漫品牌国潮文化节
$wsExhibit.Range("F14").Value = 5396  # 苏州·星部落&青铜树动漫嘉年华

# --- 演出 sheet (sheetId 2) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 103       # 苏州·乐队署同人only live Band Set二次元乐队拼盘

# --- 全部类型 sheet (sheetId 4, aggregate of all events) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 103        # 苏州·乐队署同人only live Band Set二次元乐队拼盘
$wsAll.Range("F8").Value = 131        # 苏州·无限次元夜场
$wsAll.Range("F10").Value = 141       # 苏州·漫语堂动漫嘉年华
$wsAll.Range("F11").Value = 330       # 苏州·第三届华盟国漫次元嘉年华
$wsAll.Range("F15").Value = 11578     # 苏州·COME IN JOY 动漫品牌国潮文化节
$wsAll.Range("F17").Value = 5396      # 苏州·星部落&青铜树动漫嘉年华
